# Add a new "Form Tag" column (column V) to the Investor KYC import
# template, with sample values for the first two investor rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell right after the existing "Custom field   2" column (U).
$ws.Range("V1").Value = "Form Tag"

# Sample data for the first two investors.
$ws.Range("V2").Value = "Default"
$ws.Range("V3").Value = "Gift City"

# Leave the new column blank for the remaining sample rows (3 & 4),
# matching the other optional/sparse columns in this sheet.

# Move the selection to the newly added cell, like a user would after
# typing the new column in.
$ws.Range("V2").Select()
